# Update "want to go" head-count figures (column F) and a couple of
# minimum-ticket-price values (column G) across the four sheets of the
# "广州-漫展信息" workbook, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value = 934
$ws.Range("F7").Value = 471
$ws.Range("F9").Value = 2222
$ws.Range("F10").Value = 636
$ws.Range("F13").Value = 1115
$ws.Range("F15").Value = 2242
$ws.Range("F16").Value = 691
$ws.Range("F17").Value = 13929
$ws.Range("F18").Value = 9
$ws.Range("F19").Value = 1302
$ws.Range("F20").Value = 65
$ws.Range("F21").Value = 568
$ws.Range("F22").Value = 140
$ws.Range("F23").Value = 38
$ws.Range("G24").Value = 39.9
$ws.Range("F27").Value = 282
$ws.Range("F29").Value = 4
$ws.Range("F31").Value = 32

# --- Sheet 2: 演出 (Performances) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F13").Value = 3

# --- Sheet 3: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5728
$ws.Range("F4").Value = 476

# --- Sheet 4: 全部类型 (All types, aggregate) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 476
$ws.Range("F8").Value = 934
$ws.Range("F10").Value = 471
$ws.Range("F12").Value = 2222
$ws.Range("F13").Value = 636
$ws.Range("F18").Value = 1115
$ws.Range("F23").Value = 2242
$ws.Range("F24").Value = 691
$ws.Range("F27").Value = 1302
$ws.Range("F28").Value = 65
$ws.Range("F29").Value = 568
$ws.Range("F30").Value = 140
$ws.Range("F31").Value = 38
$ws.Range("G32").Value = 39.9
$ws.Range("F38").Value = 282
$ws.Range("F49").Value = 32
